$wb = $excel.ActiveWorkbook

# ---- Sheet ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 27
$ws.Range("H27").Value = 286.66666
$ws.Range("J27").Value = 286.66666
$ws.Range("L27").Value = 859.9999799999999
$ws.Range("N27").Value = -1061.99998

# Row 43
$ws.Range("H43").Value = 4260.476
$ws.Range("I43").Value = 4293.2
$ws.Range("K43").Value = 4293.2
$ws.Range("M43").Value = -4224.2

# Row 51
$ws.Range("H51").Value = 8337633.5
$ws.Range("I51").Value = 4500
$ws.Range("J51").Value = 13893056
$ws.Range("K51").Value = 4500
$ws.Range("L51").Value = 13893056
$ws.Range("M51").Value = -4016
$ws.Range("N51").Value = -13894024

# Row 61
$ws.Range("H61").Value = 1401.4
$ws.Range("I61").Value = 1663.3334
$ws.Range("K61").Value = 4990.0002
$ws.Range("M61").Value = -4818.0002

# Row 100
$ws.Range("H100").Value = 2492.6428
$ws.Range("J100").Value = 3526
$ws.Range("L100").Value = 3526
$ws.Range("N100").Value = -4608

# Row 107
$ws.Range("H107").Value = 560.73914
$ws.Range("I107").Value = 583.5714
$ws.Range("J107").Value = 321
$ws.Range("K107").Value = 583.5714
$ws.Range("L107").Value = 321
$ws.Range("M107").Value = 1336.4286
$ws.Range("N107").Value = -4161


# ---- Sheet ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 3825.9333
$ws.Range("I61").Value = 3825.9333
$ws.Range("K61").Value = 3825.9333
$ws.Range("M61").Value = -3613.9333

# Row 74
$ws.Range("H74").Value = 1449.3846
$ws.Range("I74").Value = 1449.3846
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1449.3846
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -575.3846000000001
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 1449.3846
$ws.Range("I77").Value = 1449.3846
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 7246.923000000001
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -2878.923000000001
$ws.Range("N77").ClearContents()

# Row 110
$ws.Range("H110").Value = 982.7143
$ws.Range("I110").Value = 978
$ws.Range("K110").Value = 978
$ws.Range("M110").Value = 1067

# Row 122
$ws.Range("H122").Value = 21759.455
$ws.Range("I122").Value = 21759.455
$ws.Range("K122").Value = 65278.36500000001
$ws.Range("M122").Value = -62828.36500000001

# Row 136
$ws.Range("H136").Value = 3825.9333
$ws.Range("I136").Value = 3825.9333
$ws.Range("K136").Value = 11477.7999
$ws.Range("M136").Value = -8927.7999


# ---- Sheet BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 6540942.5
$ws.Range("I20").Value = 11114417
$ws.Range("J20").Value = 7407
$ws.Range("K20").Value = 11114417
$ws.Range("L20").Value = 7407
$ws.Range("M20").Value = -11114170
$ws.Range("N20").Value = -7901


# ---- Sheet CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2596.7827
$ws.Range("I31").Value = 1741.6364
$ws.Range("K31").Value = 1741.6364
$ws.Range("M31").Value = -1446.6364

# Row 34
$ws.Range("H34").Value = 2596.7827
$ws.Range("I34").Value = 1741.6364
$ws.Range("K34").Value = 1741.6364
$ws.Range("M34").Value = -1539.6364

# Row 58
$ws.Range("H58").Value = 5595.2915
$ws.Range("I58").Value = 6473.1055
$ws.Range("K58").Value = 6473.1055
$ws.Range("M58").Value = -6270.1055

# Row 62
$ws.Range("H62").Value = 3748.8
$ws.Range("I62").Value = 3686.25
$ws.Range("K62").Value = 3686.25
$ws.Range("M62").Value = -3062.25

# Row 65
$ws.Range("H65").Value = 3748.8
$ws.Range("I65").Value = 3686.25
$ws.Range("K65").Value = 18431.25
$ws.Range("M65").Value = -15311.25

# Row 132
$ws.Range("H132").Value = 2017.359
$ws.Range("I132").Value = 1938.8684
$ws.Range("K132").Value = 5816.6052
$ws.Range("M132").Value = -3286.6052

# Row 134
$ws.Range("H134").Value = 3004.8333
$ws.Range("I134").Value = 2881.0625
$ws.Range("K134").Value = 8643.1875
$ws.Range("M134").Value = -6108.1875

# Row 136
$ws.Range("H136").Value = 5595.2915
$ws.Range("I136").Value = 6473.1055
$ws.Range("K136").Value = 19419.3165
$ws.Range("M136").Value = -16869.3165


# ---- Sheet CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 499.85715
$ws.Range("J12").Value = 183.8
$ws.Range("L12").Value = 551.4000000000001
$ws.Range("N12").Value = -897.4000000000001

# Row 68
$ws.Range("H68").Value = 1329.2667
$ws.Range("I68").Value = 504.375
$ws.Range("J68").Value = 2272
$ws.Range("K68").Value = 1513.125
$ws.Range("L68").Value = 6816
$ws.Range("M68").Value = -702.125
$ws.Range("N68").Value = -8438

# Row 70
$ws.Range("H70").Value = 6000
$ws.Range("I70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("M70").ClearContents()

# Row 71
$ws.Range("H71").Value = 1329.2667
$ws.Range("I71").Value = 504.375
$ws.Range("J71").Value = 2272
$ws.Range("K71").Value = 4539.375
$ws.Range("L71").Value = 20448
$ws.Range("M71").Value = -483.375
$ws.Range("N71").Value = -28560

# Row 73
$ws.Range("H73").Value = 6000
$ws.Range("I73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("M73").ClearContents()

# Row 75
$ws.Range("H75").Value = 6063.6665
$ws.Range("I75").Value = 4270
$ws.Range("J75").Value = 6576.143
$ws.Range("K75").Value = 12810
$ws.Range("L75").Value = 19728.429
$ws.Range("M75").Value = -11812
$ws.Range("N75").Value = -21724.429

# Row 78
$ws.Range("H78").Value = 6063.6665
$ws.Range("I78").Value = 4270
$ws.Range("J78").Value = 6576.143
$ws.Range("K78").Value = 38430
$ws.Range("L78").Value = 59185.287
$ws.Range("M78").Value = -33438
$ws.Range("N78").Value = -69169.287

# Row 81
$ws.Range("H81").Value = 50004020
$ws.Range("I81").Value = 1058.8
$ws.Range("J81").Value = 100006984
$ws.Range("K81").Value = 3176.4
$ws.Range("L81").Value = 300020952
$ws.Range("M81").Value = -2053.4
$ws.Range("N81").Value = -300023198

# Row 84
$ws.Range("H84").Value = 50004020
$ws.Range("I84").Value = 1058.8
$ws.Range("J84").Value = 100006984
$ws.Range("K84").Value = 9529.199999999999
$ws.Range("L84").Value = 900062856
$ws.Range("M84").Value = -3913.199999999999
$ws.Range("N84").Value = -900074088

# Row 98
$ws.Range("H98").Value = 632.6667
$ws.Range("I98").Value = 500
$ws.Range("J98").Value = 659.2
$ws.Range("K98").Value = 1500
$ws.Range("L98").Value = 1977.6
$ws.Range("M98").Value = -2
$ws.Range("N98").Value = -4973.6


# ---- Sheet GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 96
$ws.Range("H96").Value = 69507.336
$ws.Range("J96").Value = 69507.336
$ws.Range("L96").Value = 69507.336
$ws.Range("N96").Value = -74999.336

# Row 122
$ws.Range("H122").Value = 1422.8889
$ws.Range("I122").Value = 1422.8889
$ws.Range("K122").Value = 4268.6667
$ws.Range("M122").Value = -1818.6667

# Row 132
$ws.Range("H132").Value = 2231.1052
$ws.Range("I132").Value = 2231.1052
$ws.Range("K132").Value = 6693.3156
$ws.Range("M132").Value = -4163.3156


# ---- Sheet LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 1094.9375
$ws.Range("I16").Value = 1208.8572
$ws.Range("J16").Value = 297.5
$ws.Range("K16").Value = 1208.8572
$ws.Range("L16").Value = 297.5
$ws.Range("M16").Value = -1038.8572
$ws.Range("N16").Value = -637.5

# Row 61
$ws.Range("H61").Value = 8935
$ws.Range("I61").Value = 9558.125
$ws.Range("K61").Value = 9558.125
$ws.Range("M61").Value = -9356.125

# Row 113
$ws.Range("H113").Value = 8935
$ws.Range("I113").Value = 9558.125
$ws.Range("K113").Value = 9558.125
$ws.Range("M113").Value = -7388.125

# Row 136
$ws.Range("H136").Value = 2838.0588
$ws.Range("I136").Value = 2126.6667
$ws.Range("K136").Value = 6380.000100000001
$ws.Range("M136").Value = -3830.000100000001


# ---- Sheet WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 132
$ws.Range("H132").Value = 3987.8254
$ws.Range("I132").Value = 4345.9214
$ws.Range("J132").Value = 2465.9167
$ws.Range("K132").Value = 13037.7642
$ws.Range("L132").Value = 7397.750100000001
$ws.Range("M132").Value = -10507.7642
$ws.Range("N132").Value = -12457.7501

# Row 136
$ws.Range("H136").Value = 1307.8125
$ws.Range("I136").Value = 513.46155
$ws.Range("K136").Value = 1540.38465
$ws.Range("M136").Value = 1009.61535

